$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, copying the formatting from the adjacent
# header cell (G1) so it reuses the existing bold/bordered header style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill H2:H8 with 0 (numeric), matching the new "Save" column values.
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
